# Weekly refresh of the "Zapallo italiano" (Vega Central Mapocho de Santiago)
# price sheet: a new week's record is inserted as the new first data row of
# the price-history block (row 383), pushing the previously-existing rows
# 383-398 down to 384-399 and growing the used range by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 383; everything below (old rows 383-398) shifts
# down to 384-399, carrying its formatting (incl. the date-style on col D).
$ws.Rows.Item(383).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(383, 1).Value2  = 9
$ws.Cells.Item(383, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(383, 3).Value2  = "Metropolitana"
$ws.Cells.Item(383, 4).Value2  = 44753
$ws.Cells.Item(383, 5).Value2  = 13
$ws.Cells.Item(383, 6).Value2  = 100112032
$ws.Cells.Item(383, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(383, 8).Value2  = "Sin especificar"
$ws.Cells.Item(383, 9).Value2  = "Primera"
$ws.Cells.Item(383, 10).Value2 = 61
$ws.Cells.Item(383, 11).Value2 = 10000
$ws.Cells.Item(383, 12).Value2 = 12000
$ws.Cells.Item(383, 13).Value2 = 11016
$ws.Cells.Item(383, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(383, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(383, 16).Value2 = 220
$ws.Cells.Item(383, 17).Value2 = 50
$ws.Cells.Item(383, 18).Value2 = "Hortaliza"
